# Add new "I0" (column I) and "IF" (column J) stat columns to the sheet,
# mirroring the existing header/body formatting used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): clone H1's formatting (bold font, thin border,
#     centered/top-aligned) onto I1 and J1, then set their text. ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Body rows 2-41: new numeric values for I (I0) and J (IF). ---
$data = @(
    @{ Row = 2;  I = 1; J = 2 },
    @{ Row = 3;  I = 1; J = 5 },
    @{ Row = 4;  I = 1; J = 6 },
    @{ Row = 5;  I = 1; J = 6 },
    @{ Row = 6;  I = 1; J = 5 },
    @{ Row = 7;  I = 1; J = 6 },
    @{ Row = 8;  I = 1; J = 5 },
    @{ Row = 9;  I = 1; J = 5 },
    @{ Row = 10; I = 1; J = 5 },
    @{ Row = 11; I = 1; J = 1 },
    @{ Row = 12; I = 1; J = 5 },
    @{ Row = 13; I = 1; J = 6 },
    @{ Row = 14; I = 1; J = 5 },
    @{ Row = 15; I = 1; J = 5 },
    @{ Row = 16; I = 2; J = 4 },
    @{ Row = 17; I = 1; J = 4 },
    @{ Row = 18; I = 1; J = 5 },
    @{ Row = 19; I = 1; J = 5 },
    @{ Row = 20; I = 4; J = 5 },
    @{ Row = 21; I = 9; J = 9 },
    @{ Row = 22; I = 1; J = 5 },
    @{ Row = 23; I = 1; J = 4 },
    @{ Row = 24; I = 8; J = 8 },
    @{ Row = 25; I = 4; J = 5 },
    @{ Row = 26; I = 8; J = 8 },
    @{ Row = 27; I = 8; J = 8 },
    @{ Row = 28; I = 7; J = 7 },
    @{ Row = 29; I = 8; J = 8 },
    @{ Row = 30; I = 8; J = 9 },
    @{ Row = 31; I = 8; J = 9 },
    @{ Row = 32; I = 8; J = 9 },
    @{ Row = 33; I = 9; J = 9 },
    @{ Row = 34; I = 8; J = 8 },
    @{ Row = 35; I = 9; J = 9 },
    @{ Row = 36; I = 1; J = 5 },
    @{ Row = 37; I = 7; J = 8 },
    @{ Row = 38; I = 1; J = 2 },
    @{ Row = 39; I = 7; J = 7 },
    @{ Row = 40; I = 6; J = 7 },
    @{ Row = 41; I = 7; J = 7 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I   # column I
    $ws.Cells.Item($r, 10).Value = $entry.J  # column J
}
